# Apply updated crypto price/volume figures to the "cryptos" sheet.
# Price (column D) values that parse as plain numbers must be written as
# explicit text (Excel auto-converts numeric-looking strings otherwise,
# which would both reclassify the cell type and drop trailing zeros such
# as "572.50" -> 572.5). We flip the cell to Text format just long enough
# to assign the literal string, then clear the format again so the cell
# keeps its original (default) style - only the stored value changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "63.423.57"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "2.479.86"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  +0.21%  "
Set-TextValue "D5" "572.50"
$ws.Range("E5").Value = "  +1.65%  "
Set-TextValue "D6" "149.37"
$ws.Range("E6").Value = "  +4.29%  "
$ws.Range("E7").Value = "  -0.14%  "
Set-TextValue "D8" "0.539"
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("E9").Value = "  +4.18%  "
Set-TextValue "D11" "0.364"
$ws.Range("E11").Value = "  +4.15%  "
Set-TextValue "D12" "5.34"
$ws.Range("E12").Value = "  +2.71%  "
Set-TextValue "D13" "27.19"
$ws.Range("E13").Value = "  +5.50%  "
Set-TextValue "D14" "0.0000184"
$ws.Range("E14").Value = "  +6.32%  "
$ws.Range("D15").Value = "2.920.74"
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").Value = "63.624.99"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("D17").Value = "2.508.87"
$ws.Range("E17").Value = "  +2.64%  "
Set-TextValue "D18" "11.54"
$ws.Range("E18").Value = "  +2.35%  "
Set-TextValue "D19" "7.24"
$ws.Range("E19").Value = "  +5.88%  "
Set-TextValue "D20" "4.24"
$ws.Range("E20").Value = "  +2.62%  "
Set-TextValue "D21" "328.50"
$ws.Range("E21").Value = "  +1.61%  "
Set-TextValue "D22" "0.998"
$ws.Range("E22").Value = "  -0.20%  "
Set-TextValue "D23" "1.88"
$ws.Range("E23").Value = "  +7.57%  "
Set-TextValue "D24" "67.39"
$ws.Range("E24").Value = "  +1.25%  "
Set-TextValue "D25" "641.64"
$ws.Range("E25").Value = "  +14.24%  "
Set-TextValue "D26" "0.0000105"
$ws.Range("E26").Value = "  +11.68%  "
Set-TextValue "D27" "8.75"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").Value = "2.608.53"
$ws.Range("E28").Value = "  +2.84%  "
Set-TextValue "D29" "1.52"
$ws.Range("E29").Value = "  +9.32%  "
Set-TextValue "D30" "8.53"
$ws.Range("E30").Value = "  +4.30%  "
Set-TextValue "D31" "0.997"
$ws.Range("E31").Value = "  -0.14%  "
Set-TextValue "D32" "0.144"
$ws.Range("E32").Value = "  -2.16%  "
Set-TextValue "D33" "1.92"
$ws.Range("E33").Value = "  +3.02%  "
Set-TextValue "D34" "5.17"
$ws.Range("E34").Value = "  +8.57%  "
Set-TextValue "D35" "1.54"
$ws.Range("E35").Value = "  +2.81%  "
$ws.Range("E36").Value = "  -0.25%  "
Set-TextValue "D37" "0.387"
$ws.Range("E37").Value = "  +2.13%  "
Set-TextValue "D38" "5.50"
$ws.Range("E38").Value = "  +1.08%  "
Set-TextValue "D39" "18.93"
$ws.Range("E39").Value = "  +2.30%  "
Set-TextValue "D40" "1.85"
$ws.Range("E40").Value = "  +1.81%  "
Set-TextValue "D41" "147.42"
$ws.Range("E41").Value = "  -4.31%  "
Set-TextValue "D42" "2.66"
$ws.Range("E42").Value = "  +17.59%  "
$ws.Range("E43").Value = "  +0.25%  "
Set-TextValue "D44" "151.99"
$ws.Range("E44").Value = "  +2.99%  "
Set-TextValue "D45" "3.77"
$ws.Range("E45").Value = "  +3.92%  "
Set-TextValue "D46" "0.0552"
$ws.Range("E46").Value = "  +4.60%  "
Set-TextValue "D47" "21.06"
$ws.Range("E47").Value = "  +6.15%  "
Set-TextValue "D48" "0.612"
$ws.Range("E48").Value = "  +3.32%  "
Set-TextValue "D49" "0.0238"
$ws.Range("E49").Value = "  +5.12%  "
Set-TextValue "D50" "0.0928"
$ws.Range("E50").Value = "  +0.67%  "
Set-TextValue "D51" "0.743"
$ws.Range("E51").Value = "  +4.59%  "
